# Commit: "Changed min order on cheap components to 10"
#
# The Bill-of-Materials sheet has a "Quantity" (order quantity) column D.
# A number of cheap / commodity passive components (capacitors, resistors,
# a transistor, etc.) had a minimum-order quantity of 3 (one row had 6);
# bump all of those specific rows' minimum order quantity up to 10.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Bill of Materials-S6C")

# Rows whose Quantity (column D) changes to 10.
$rowsToUpdate = @(18, 43, 45, 47, 48, 50, 53, 54, 56, 57, 59, 61, 62, 63, 64, 66, 67, 68, 69)

foreach ($row in $rowsToUpdate) {
    $ws.Range("D$row").Value = 10
}

# The merged header/section cells got re-serialized in a different (but
# equivalent) order. Re-create each merge so the on-disk order matches.
$mergedRanges = @(
    "A2:H2",
    "A13:H13",
    "A15:H15",
    "A19:H19",
    "A21:H21",
    "A52:H52",
    "A24:H24",
    "A32:H32",
    "A36:H36",
    "A38:H38",
    "A40:H40"
)

foreach ($r in $mergedRanges) {
    $ws.Range($r).UnMerge()
}
foreach ($r in $mergedRanges) {
    $ws.Range($r).Merge()
}
